$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '305.01'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2.45%'
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '6.35%'
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.144'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '2.59%'
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07837'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.95%'
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('B6').NumberFormat = 'General'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('C6').NumberFormat = 'General'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.618'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2.98%'
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('B7').NumberFormat = 'General'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C7').NumberFormat = 'General'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.056'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '13.53%'
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('B8').NumberFormat = 'General'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('C8').NumberFormat = 'General'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.1307'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '8.61%'
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('B9').NumberFormat = 'General'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('C9').NumberFormat = 'General'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1875'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.96%'
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('B10').NumberFormat = 'General'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('C10').NumberFormat = 'General'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09181'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '3.71%'
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('B11').NumberFormat = 'General'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('C11').NumberFormat = 'General'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.04160'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.02%'
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('B12').NumberFormat = 'General'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C12').NumberFormat = 'General'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1044'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.92%'
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('B13').NumberFormat = 'General'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C13').NumberFormat = 'General'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001291'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '1.47%'
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'TigerCash'
$ws.Range('B14').NumberFormat = 'General'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C14').NumberFormat = 'General'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.005820'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.33%'
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'UpBots'
$ws.Range('B15').NumberFormat = 'General'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('C15').NumberFormat = 'General'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.007445'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1,899.16%'
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('B16').NumberFormat = 'General'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C16').NumberFormat = 'General'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.368'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.72%'
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('B17').NumberFormat = 'General'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('C17').NumberFormat = 'General'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.433'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.52%'
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3379'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '1.28%'
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.046'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.75%'
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1380'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-2.57%'
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2806'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-6.28%'
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04185'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '3.55%'
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001275'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.06%'
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '14.81%'
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001342'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '9.23%'
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02585'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '6.63%'
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05351'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '2.54%'
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.005597'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-4.84%'
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007772'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.37%'
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1381'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '3.62%'
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007305'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-0.71%'
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '6.49%'
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3023'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '1.58%'
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006692'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '5.86%'
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.46%'
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06572'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '45.71%'
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.003974'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-5.20%'
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002087'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.46%'
$ws.Range('E50').NumberFormat = 'General'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.46%'
$ws.Range('E51').NumberFormat = 'General'
